$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Change discount rates (column L, rows 2-7) from 2.5% to 5.0%
$ws.Range("L2:L7").Value = 0.05

# Update the view state: scroll position (topLeftCell B1 -> C1) and active cell selection (R11 -> L8)
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("L8").Select() | Out-Null
